$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 2; $r -le 78; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $current = $cell.Value2()
    $cell.Value2 = $current + 2
}
